$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values are stored as text in the source data (e.g. "30.315.80",
# "0.9999"). A plain .Value assignment lets Excel auto-coerce numeric-looking
# strings into numbers, so force text via NumberFormat="@" then restore the default
# "Normal" style (keeps styles.xml untouched / no stray style index on the cell).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.315.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4827'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2891'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06593'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.878.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.93'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07393'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.180'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6600'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.284.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007731'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.472'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.136.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '195.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.415'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.927'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.437'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.270'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09134'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.048'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05051'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7412'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.140'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.706'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01837'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.634'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9143'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.075'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.53'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.877'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4321'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.629'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1348'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.569'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.900'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05716'
$ws.Range("D51").Style = "Normal"

# Volume(1h) column (E): values like "  -1.22%  " are never numeric-looking, so a
# plain assignment already round-trips as text.
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -2.14%  '
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("E25").Value = '  -2.49%  '
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("E27").Value = '  -3.77%  '
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("E29").Value = '  -2.95%  '
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("E33").Value = '  -4.53%  '
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("E35").Value = '  +1.88%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  -2.94%  '
$ws.Range("E47").Value = '  +8.83%  '
$ws.Range("E48").Value = '  -13.27%  '
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("E50").Value = '  -5.35%  '
$ws.Range("E51").Value = '  -2.76%  '
